$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '93.224.70'
$ws.Cells.Item(2, 5).Value = '  +1.77%  '

$ws.Cells.Item(3, 4).Value = '3.129.39'
$ws.Cells.Item(3, 5).Value = '  -1.91%  '

$ws.Cells.Item(4, 5).Value = '  -0.02%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '243.84'
$ws.Cells.Item(5, 5).Value = '  +1.45%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '617.39'
$ws.Cells.Item(6, 5).Value = '  -0.76%  '

$ws.Cells.Item(7, 5).Value = '  -2.37%  '

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.405'
$ws.Cells.Item(8, 5).Value = '  +8.06%  '

$ws.Cells.Item(9, 5).Value = '  -0.03%  '

$ws.Cells.Item(10, 4).Value = '3.124.76'
$ws.Cells.Item(10, 5).Value = '  -1.87%  '

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.736'
$ws.Cells.Item(11, 5).Value = '  -2.15%  '

$ws.Cells.Item(12, 5).Value = '  -1.30%  '

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '0.0000255'
$ws.Cells.Item(13, 5).Value = '  +2.45%  '

$ws.Cells.Item(14, 4).Value = '92.795.80'
$ws.Cells.Item(14, 5).Value = '  +1.61%  '

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '34.66'
$ws.Cells.Item(15, 5).Value = '  -2.99%  '

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '5.51'
$ws.Cells.Item(16, 5).Value = '  -0.67%  '

$ws.Cells.Item(17, 4).Value = '3.701.43'
$ws.Cells.Item(17, 5).Value = '  -1.69%  '

$ws.Cells.Item(18, 4).Value = '3.105.27'
$ws.Cells.Item(18, 5).Value = '  -1.89%  '

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '3.79'
$ws.Cells.Item(19, 5).Value = '  +1.36%  '

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '14.81'
$ws.Cells.Item(20, 5).Value = '  -3.37%  '

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '5.83'
$ws.Cells.Item(21, 5).Value = '  -2.20%  '

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '0.0000207'
$ws.Cells.Item(22, 5).Value = '  -0.27%  '

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '9.48'
$ws.Cells.Item(23, 5).Value = '  +2.54%  '

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '450.56'
$ws.Cells.Item(24, 5).Value = '  +0.92%  '

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '5.84'
$ws.Cells.Item(25, 5).Value = '  -4.99%  '

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '87.27'
$ws.Cells.Item(26, 5).Value = '  -2.64%  '

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '11.85'
$ws.Cells.Item(27, 5).Value = '  -2.95%  '

$ws.Cells.Item(28, 4).Value = '3.283.85'
$ws.Cells.Item(28, 5).Value = '  -1.67%  '

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '1.00'
$ws.Cells.Item(29, 5).Value = '  +0.32%  '

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '0.137'
$ws.Cells.Item(30, 5).Value = '  +3.97%  '

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '0.229'
$ws.Cells.Item(31, 5).Value = '  -1.66%  '

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '0.169'
$ws.Cells.Item(32, 5).Value = '  -1.17%  '

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '9.27'
$ws.Cells.Item(33, 5).Value = '  -1.85%  '

$ws.Cells.Item(34, 5).Value = '  +10.59%  '

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '8.06'
$ws.Cells.Item(35, 5).Value = '  +2.97%  '

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.165'
$ws.Cells.Item(36, 5).Value = '  -3.14%  '

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '26.32'
$ws.Cells.Item(37, 5).Value = '  -1.15%  '

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '4.20'
$ws.Cells.Item(38, 5).Value = '  +10.23%  '

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '1.92'

$ws.Cells.Item(40, 2).Value = 'Bittensor'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '487.21'
$ws.Cells.Item(40, 5).Value = '  -5.37%  '

$ws.Cells.Item(41, 2).Value = 'Fetch.AI'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '1.32'
$ws.Cells.Item(41, 5).Value = '  -4.15%  '

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '3.52'
$ws.Cells.Item(42, 5).Value = '  +1.61%  '

$ws.Cells.Item(43, 5).Value = '  -4.11%  '

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '23.07'
$ws.Cells.Item(44, 5).Value = '  +4.07%  '

$ws.Cells.Item(45, 5).Value = '  -0.02%  '

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '161.33'
$ws.Cells.Item(46, 5).Value = '  +2.19%  '

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '1.94'
$ws.Cells.Item(47, 5).Value = '  +0.09%  '

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '0.696'
$ws.Cells.Item(48, 5).Value = '  -5.12%  '

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '1.39'
$ws.Cells.Item(49, 5).Value = '  -0.93%  '

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '0.0339'
$ws.Cells.Item(50, 5).Value = '  +5.55%  '

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '4.48'
$ws.Cells.Item(51, 5).Value = '  +0.00%  '
